$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2 through 18: change date serial 45182 -> 45184
$ws.Range("C2:C18").Value = 45184
